$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "{'Cus_Nme', 'Acc_num', 'Loan_Amt', 'Cust_Addr'}{'Acc_num', 'Debit_pin', 'From_AcctNum', 'To_AcctNum', 'Amt_avail'}"
$ws.Range("D3").Value = "{'Cus_Nme', 'Acc_num', 'Debit_pin', 'Bill_type', 'Max_limit'}{'Cus_Nme', 'Acc_num', 'Loan_Amt', 'Cust_Addr'}"
$ws.Range("D4").Value = "{'Cus_Nme', 'From_AcctNum', 'Amt_trnsfr', 'To_AcctNum', 'Amt_avail'}{'Cus_Nme', 'Acc_num', 'Loan_Amt', 'Cust_Addr'}"
$ws.Range("D5").Value = "{'Acc_num', 'Bill_type', 'Debit_pin', 'Amt_avail', 'From_AcctNum'}"
$ws.Range("D6").Value = "{'Cus_Nme', 'Acc_num', 'Debit_pin', 'To_AcctNum', 'From_AcctNum'}{'Loan_Amt', 'Loan_purp', 'Cred_Score'}"
$ws.Range("D7").Value = "{'Cus_Nme', 'Acc_num', 'Loan_Amt', 'Cust_Addr'}"
$ws.Range("D8").Value = "{'Cus_Nme', 'Acc_num', 'Acc_type', 'Max_limit'}{'Cus_Nme', 'Acc_num', 'Debit_pin', 'Acc_type', 'Amt_deposit'}"
$ws.Range("D9").Value = "{'Cus_Nme', 'Acc_num', 'Debit_pin', 'Acc_type', 'Amt_deposit'}"
$ws.Range("D10").Value = "{'Acc_num'}{'Cus_Nme', 'Acc_num', 'Loan_Amt', 'Cust_Addr'}"
$ws.Range("D11").Value = "{'Loan_Amt', 'Loan_purp', 'Cred_Score'}{'Acc_num'},{'Cus_Nme', 'Acc_num', 'Debit_pin', 'To_AcctNum', 'From_AcctNum'}{'Loan_Amt', 'Loan_purp', 'Cred_Score'}"
$ws.Range("D12").Value = ""
